$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 15 (pushes the "Phase2" header block and
# everything below it down by one row).
$ws.Rows.Item(15).Insert()

# The old "users should be able to email a recipe" row (previously row 17,
# now shifted down to row 18) is no longer needed as its own line item -
# remove it entirely, which shifts rows 19+ back up by one.
$ws.Rows.Item(18).Delete()

# Populate the newly inserted row 15 with the new enhancement entry.
$ws.Cells.Item(15, 2).Value = "dishes/_form"
$ws.Cells.Item(15, 4).Value = "capability to email the recipe"
$ws.Cells.Item(15, 4).WrapText = $true

# Update the selected/active cell to reflect where the user ended up.
$ws.Range("E15").Select()
